$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from its old spot (an empty paragraph near
#    the end of the document) into the middle of the sentence "...to frame
#    the debate is based..." -- i.e. right after "frame the " and before
#    "debate". This reflects the last place the author actually edited text
#    during the Feb 18 2019 lecture revision.
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("skeptics of the program to frame the ")
if (-not $found) {
    throw "Could not find anchor text for bookmark relocation"
}
$r.Collapse(0)  # wdCollapseEnd
$d.Bookmarks.Add("_GoBack", $r)

# ---------------------------------------------------------------------------
# 2) Insert a new blank paragraph (matching the spacing of its neighbors)
#    between the "One reason the policy narrative of proponents..." paragraph
#    and the bold "In two or three succinct sentences..." prompt paragraph.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("and are undeserving of any inherent bias in the system that favors them over small businesses.")
if (-not $found2) {
    throw "Could not find anchor text for new blank paragraph"
}
$para = $r2.Paragraphs(1)
$para.Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 3) Collapse the "Response to post by James Tillis about policy narratives
#    in the debate about urban agriculture zones:" paragraph's several runs
#    into a single run (the text itself is unchanged).
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute(
    "Response to post by James Tillis about policy narratives in the debate about urban agriculture zones:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Response to post by James Tillis about policy narratives in the debate about urban agriculture zones:",
    2
) | Out-Null

Write-Output "done"
